# Auto-generated edit script: updates market-price-derived cells
# (currentAveragePrice / LevePrice / LeveProfit columns) across all
# leve-profit worksheets, per the scheduled market-data refresh.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 4227.161
$ws.Range("J17").Value = 4227.161
$ws.Range("L17").Value = 12681.483
$ws.Range("N17").Value = -13017.483
$ws.Range("H40").Value = 2136.889
$ws.Range("I40").Value = 2068.4443
$ws.Range("J40").Value = 2273.7778
$ws.Range("K40").Value = 2068.4443
$ws.Range("L40").Value = 2273.7778
$ws.Range("M40").Value = -1893.4443
$ws.Range("N40").Value = -2623.7778
$ws.Range("H64").Value = 4733
$ws.Range("I64").Value = 3199.5
$ws.Range("J64").Value = 5499.75
$ws.Range("K64").Value = 3199.5
$ws.Range("L64").Value = 5499.75
$ws.Range("M64").Value = -2951.5
$ws.Range("N64").Value = -5995.75
$ws.Range("H67").Value = 4733
$ws.Range("I67").Value = 3199.5
$ws.Range("J67").Value = 5499.75
$ws.Range("K67").Value = 3199.5
$ws.Range("L67").Value = 5499.75
$ws.Range("M67").Value = -2341.5
$ws.Range("N67").Value = -7215.75
$ws.Range("H80").Value = 6890.8335
$ws.Range("J80").Value = 8329.75
$ws.Range("L80").Value = 24989.25
$ws.Range("N80").Value = -26985.25
$ws.Range("H83").Value = 6890.8335
$ws.Range("J83").Value = 8329.75
$ws.Range("L83").Value = 74967.75
$ws.Range("N83").Value = -84951.75
$ws.Range("H127").Value = 5958.3335
$ws.Range("J127").Value = 4438.25
$ws.Range("L127").Value = 13314.75
$ws.Range("N127").Value = -23234.75

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 398
$ws.Range("I2").Value = 397.8
$ws.Range("K2").Value = 397.8
$ws.Range("M2").Value = -284.8
$ws.Range("H32").Value = 11738.083
$ws.Range("J32").Value = 7662.6665
$ws.Range("L32").Value = 7662.6665
$ws.Range("N32").Value = -8236.666499999999
$ws.Range("H45").Value = 4526.1113
$ws.Range("I45").Value = 3840.125
$ws.Range("K45").Value = 3840.125
$ws.Range("M45").Value = -3463.125
$ws.Range("H97").Value = 405.66666
$ws.Range("I97").Value = 405.66666
$ws.Range("K97").Value = 405.66666
$ws.Range("M97").Value = 90.33334000000002
$ws.Range("H116").Value = 398
$ws.Range("I116").Value = 397.8
$ws.Range("K116").Value = 397.8
$ws.Range("M116").Value = 1896.2
$ws.Range("H132").Value = 3749.6667
$ws.Range("I132").Value = 3624.5
$ws.Range("K132").Value = 10873.5
$ws.Range("M132").Value = -8343.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H2").Value = 63999.5
$ws.Range("J2").Value = 63999.5
$ws.Range("L2").Value = 63999.5
$ws.Range("N2").Value = -64225.5
$ws.Range("H3").Value = 398
$ws.Range("I3").Value = 397.8
$ws.Range("K3").Value = 397.8
$ws.Range("M3").Value = -283.8
$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()
$ws.Range("H49").Value = 0
$ws.Range("I49").Value = 0
$ws.Range("K49").Value = 0
$ws.Range("M49").ClearContents()
$ws.Range("H80").Value = 269.05884
$ws.Range("I80").Value = 346.33334
$ws.Range("J80").Value = 226.90909
$ws.Range("K80").Value = 346.33334
$ws.Range("L80").Value = 226.90909
$ws.Range("M80").Value = 651.66666
$ws.Range("N80").Value = -2222.90909
$ws.Range("H83").Value = 269.05884
$ws.Range("I83").Value = 346.33334
$ws.Range("J83").Value = 226.90909
$ws.Range("K83").Value = 1731.6667
$ws.Range("L83").Value = 1134.54545
$ws.Range("M83").Value = 3260.3333
$ws.Range("N83").Value = -11118.54545
$ws.Range("H133").Value = 62332
$ws.Range("I133").Value = 59498.5
$ws.Range("J133").Value = 67999
$ws.Range("K133").Value = 59498.5
$ws.Range("L133").Value = 67999
$ws.Range("M133").Value = -54438.5
$ws.Range("N133").Value = -78119

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1250.3334
$ws.Range("I16").Value = 957.1429000000001
$ws.Range("J16").Value = 2276.5
$ws.Range("K16").Value = 957.1429000000001
$ws.Range("L16").Value = 2276.5
$ws.Range("M16").Value = -670.1429000000001
$ws.Range("N16").Value = -2850.5
$ws.Range("H31").Value = 2332.3635
$ws.Range("I31").Value = 1919.6666
$ws.Range("K31").Value = 1919.6666
$ws.Range("M31").Value = -1624.6666
$ws.Range("H34").Value = 2332.3635
$ws.Range("I34").Value = 1919.6666
$ws.Range("K34").Value = 1919.6666
$ws.Range("M34").Value = -1717.6666
$ws.Range("H105").Value = 3783.8262
$ws.Range("I105").Value = 3155.2222
$ws.Range("J105").Value = 4187.9287
$ws.Range("K105").Value = 3155.2222
$ws.Range("L105").Value = 4187.9287
$ws.Range("M105").Value = -1408.2222
$ws.Range("N105").Value = -7681.9287
$ws.Range("H108").Value = 48999.5
$ws.Range("J108").Value = 48999.5
$ws.Range("L108").Value = 48999.5
$ws.Range("N108").Value = -56679.5
$ws.Range("H113").Value = 1250.3334
$ws.Range("I113").Value = 957.1429000000001
$ws.Range("J113").Value = 2276.5
$ws.Range("K113").Value = 957.1429000000001
$ws.Range("L113").Value = 2276.5
$ws.Range("M113").Value = 1212.8571
$ws.Range("N113").Value = -6616.5

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 79982.336
$ws.Range("J37").Value = 79982.336
$ws.Range("L37").Value = 239947.008
$ws.Range("N37").Value = -240171.008
$ws.Range("H109").Value = 2877.5557
$ws.Range("I109").Value = 3056.125
$ws.Range("J109").Value = 1449
$ws.Range("K109").Value = 9168.375
$ws.Range("L109").Value = 4347
$ws.Range("M109").Value = -8128.375
$ws.Range("N109").Value = -6427
$ws.Range("H114").Value = 2023.8823
$ws.Range("J114").Value = 1947.6
$ws.Range("L114").Value = 5842.799999999999
$ws.Range("N114").Value = -12350.8
$ws.Range("H117").Value = 13080.625
$ws.Range("I117").Value = 551.2
$ws.Range("K117").Value = 1653.6
$ws.Range("M117").Value = 1788.4
$ws.Range("H129").Value = 2849.3333
$ws.Range("I129").Value = 1230
$ws.Range("J129").Value = 3659
$ws.Range("K129").Value = 3690
$ws.Range("L129").Value = 10977
$ws.Range("M129").Value = 1310
$ws.Range("N129").Value = -20977
$ws.Range("H130").Value = 2403
$ws.Range("J130").Value = 2978.6667
$ws.Range("L130").Value = 8936.000100000001
$ws.Range("N130").Value = -18976.0001
$ws.Range("H131").Value = 2879.739
$ws.Range("J131").Value = 2879.739
$ws.Range("L131").Value = 8639.217000000001
$ws.Range("N131").Value = -18719.217
$ws.Range("H133").Value = 12997.857
$ws.Range("I133").Value = 3661.6667
$ws.Range("K133").Value = 10985.0001
$ws.Range("M133").Value = -5925.000100000001
$ws.Range("H137").Value = 2345.6155
$ws.Range("I137").Value = 1610.6666
$ws.Range("K137").Value = 4831.9998
$ws.Range("M137").Value = 268.0002000000004

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 951.6667
$ws.Range("I80").Value = 777.5
$ws.Range("J80").Value = 1300
$ws.Range("K80").Value = 777.5
$ws.Range("L80").Value = 1300
$ws.Range("M80").Value = 220.5
$ws.Range("N80").Value = -3296
$ws.Range("H83").Value = 951.6667
$ws.Range("I83").Value = 777.5
$ws.Range("J83").Value = 1300
$ws.Range("K83").Value = 3887.5
$ws.Range("L83").Value = 6500
$ws.Range("M83").Value = 1104.5
$ws.Range("N83").Value = -16484
$ws.Range("H113").Value = 1052.8
$ws.Range("I113").Value = 1052.8
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1052.8
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 1117.2
$ws.Range("N113").ClearContents()
$ws.Range("H123").Value = 149999.5
$ws.Range("J123").Value = 149999.5
$ws.Range("L123").Value = 149999.5
$ws.Range("N123").Value = -154899.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1798.8572
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("H27").Value = 1798.8572
$ws.Range("I27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("M27").ClearContents()
$ws.Range("H46").Value = 33723.25
$ws.Range("I46").Value = 64446.625
$ws.Range("K46").Value = 64446.625
$ws.Range("M46").Value = -64258.625
$ws.Range("H82").Value = 1237.5
$ws.Range("I82").Value = 1345.375
$ws.Range("K82").Value = 1345.375
$ws.Range("M82").Value = -984.375
$ws.Range("H85").Value = 1237.5
$ws.Range("I85").Value = 1345.375
$ws.Range("K85").Value = 1345.375
$ws.Range("M85").Value = -97.375
$ws.Range("H136").Value = 2899.4
$ws.Range("I136").Value = 2468.0908
$ws.Range("K136").Value = 7404.2724
$ws.Range("M136").Value = -4854.2724
$ws.Range("H137").Value = 97693.5
$ws.Range("I137").Value = 85000
$ws.Range("J137").Value = 110387
$ws.Range("K137").Value = 85000
$ws.Range("L137").Value = 110387
$ws.Range("M137").Value = -79900
$ws.Range("N137").Value = -120587

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 18096.385
$ws.Range("I62").Value = 33199
$ws.Range("J62").Value = 13565.6
$ws.Range("K62").Value = 33199
$ws.Range("L62").Value = 13565.6
$ws.Range("M62").Value = -32575
$ws.Range("N62").Value = -14813.6
$ws.Range("H65").Value = 18096.385
$ws.Range("I65").Value = 33199
$ws.Range("J65").Value = 13565.6
$ws.Range("K65").Value = 165995
$ws.Range("L65").Value = 67828
$ws.Range("M65").Value = -162875
$ws.Range("N65").Value = -74068
$ws.Range("H107").Value = 2000
$ws.Range("I107").Value = 2000
$ws.Range("K107").Value = 6000
$ws.Range("M107").Value = -4080
$ws.Range("H122").Value = 14167.333
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("H138").Value = 89998
$ws.Range("J138").Value = 89998
$ws.Range("L138").Value = 89998
$ws.Range("N138").Value = -100278
